$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: add H21 label
$ws.Range("H21").Value = "(Design) 2.2 data viz"

# Row 22: add E22 value (drives F22/G22 recalculation), add H22 label, update J22
$ws.Range("E22").Value = 57
$ws.Range("H22").Value = "2.3 GUIs"
$ws.Range("J22").Value = 57

# Row 23: update J23 formula (base for the shared formula chain J24:J37)
$ws.Range("J23").Formula = "=J22+4"

# Update selection/view state
$ws.Range("CF22").Select()
